# Apply cryptocurrency price/volume updates from the Jul 6 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold numeric-looking text; force text entry
# via a temporary "@" number format, then restore the default "Normal" style so the
# cell formatting/style index is left exactly as it was (matches source which had no
# explicit style on these cells).
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "56.233.01"
Set-TextValue "E2" "  +3.45%  "
Set-TextValue "D3" "2.969.35"
Set-TextValue "E3" "  +2.36%  "
Set-TextValue "E4" "  +0.09%  "
Set-TextValue "D5" "501.54"
Set-TextValue "E5" "  +5.91%  "
Set-TextValue "D6" "134.25"
Set-TextValue "E6" "  +6.19%  "
Set-TextValue "E7" "  +0.01%  "
Set-TextValue "E8" "  +6.32%  "
Set-TextValue "E9" "  +10.32%  "
Set-TextValue "D10" "0.106"
Set-TextValue "E10" "  +10.30%  "
Set-TextValue "E11" "  +5.74%  "
Set-TextValue "E12" "  +3.11%  "
Set-TextValue "D13" "3.481.34"
Set-TextValue "E13" "  +2.62%  "
Set-TextValue "D14" "25.18"
Set-TextValue "E14" "  +9.66%  "
Set-TextValue "E15" "  +12.63%  "
Set-TextValue "D16" "56.249.96"
Set-TextValue "E16" "  +3.39%  "
Set-TextValue "D17" "2.971.03"
Set-TextValue "E17" "  +2.66%  "
Set-TextValue "D18" "5.65"
Set-TextValue "E18" "  +10.18%  "
Set-TextValue "D19" "12.29"
Set-TextValue "E19" "  +6.87%  "
Set-TextValue "D20" "7.73"
Set-TextValue "E20" "  +8.47%  "
Set-TextValue "D21" "321.62"
Set-TextValue "E21" "  +4.59%  "
Set-TextValue "D22" "0.998"
Set-TextValue "E22" "  -0.23%  "
Set-TextValue "E23" "  +5.04%  "
Set-TextValue "D24" "61.86"
Set-TextValue "E24" "  +4.77%  "
Set-TextValue "D25" "1.00"
Set-TextValue "E25" "  -0.18%  "
Set-TextValue "E26" "  +5.49%  "
Set-TextValue "D27" "0.0₃0883"
Set-TextValue "E27" "  +6.63%  "
Set-TextValue "D28" "6.47"
Set-TextValue "E28" "  +5.75%  "
Set-TextValue "D29" "6.76"
Set-TextValue "E29" "  +10.43%  "
Set-TextValue "D30" "1.18"
Set-TextValue "E30" "  +2.84%  "
Set-TextValue "D31" "1.74"
Set-TextValue "E31" "  +9.18%  "
Set-TextValue "D32" "20.37"
Set-TextValue "E32" "  +6.31%  "
Set-TextValue "D33" "158.25"
Set-TextValue "E33" "  +11.64%  "
Set-TextValue "D34" "4.44"
Set-TextValue "E34" "  +5.21%  "
Set-TextValue "D35" "1.25"
Set-TextValue "E35" "  +3.30%  "
Set-TextValue "D36" "5.52"
Set-TextValue "E36" "  +1.07%  "
Set-TextValue "E37" "  +8.58%  "
Set-TextValue "D38" "23.03"
Set-TextValue "E38" "  +2.36%  "
Set-TextValue "D39" "3.005.42"
Set-TextValue "E39" "  +2.72%  "
Set-TextValue "E40" "  +0.17%  "
Set-TextValue "D41" "36.16"
Set-TextValue "E41" "  +4.37%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D42" "0.640"
Set-TextValue "E42" "  +7.01%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D43" "2.243.56"
Set-TextValue "E43" "  +8.74%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D44" "1.39"
Set-TextValue "E44" "  +6.64%  "
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue "D45" "0.979"
Set-TextValue "E45" "  +2.03%  "
Set-TextValue "D46" "3.54"
Set-TextValue "E46" "  +3.40%  "
Set-TextValue "E47" "  +20.92%  "
Set-TextValue "D48" "5.74"
Set-TextValue "E48" "  +8.40%  "
Set-TextValue "D49" "0.0233"
Set-TextValue "E49" "  +10.18%  "
Set-TextValue "E50" "  +4.68%  "
Set-TextValue "E51" "  +8.72%  "
